$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the payment/claim records with the new user, password and claim numbers
$ws.Range("C2").Value = "ocerutti"
$ws.Range("D2").Value = "silverarrow"

$ws.Range("F3").Value = "'1220170301357"
$ws.Range("F2").Value = "'1220194200610"

# Move the active selection to F12, matching where work left off
$ws.Range("F12").Select() | Out-Null
